$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# The title text is split across three runs: "Below", " ", "section-level".
# Drop the middle (space-only) run first so the remaining "Below" and
# "section-level" runs become adjacent...
$sp = $tr.Characters(6, 1)
$sp.Text = ""

# ...then rewrite the whole (now-adjacent) range so it collapses into a
# single run containing "Below section-level".
$whole = $tr.Characters(1, $tr.Length)
$whole.Text = "Below section-level"
